$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# REPORTDATE
$ws.Range("H2").Value = "2020-09-30 00:00:00"

# BASIC_EPS
$ws.Range("I2").Value = 0.003

# TOTAL_OPERATE_INCOME
$ws.Range("K2").Value = 674062.13

# PARENT_NETPROFIT
$ws.Range("L2").Value = 100298.01

# YSTZ
$ws.Range("N2").Value = 135.4944608179

# SJLTZ
$ws.Range("O2").Value = 101.5470680587

# BPS
$ws.Range("P2").Value = 0.82515509

# MGJYXJJE
$ws.Range("Q2").Value = -0.000561950968

# XSMLL
$ws.Range("R2").Value = 72.256815555

# ISNEW (must stay text "1", not get auto-converted to the number 1)
$ws.Range("AB2").NumberFormat = "@"
$ws.Range("AB2").Value = "1"
$ws.Range("AB2").ClearFormats()

# QDATE
$ws.Range("AC2").Value = "2020Q3"

# DATATYPE
$ws.Range("AD2").Value = "2020年 三季报"

# DATAYEAR (must stay text "2020", not get auto-converted to the number 2020)
$ws.Range("AE2").NumberFormat = "@"
$ws.Range("AE2").Value = "2020"
$ws.Range("AE2").ClearFormats()
